$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 9: Erik Meurrens / Garage networking ---
# Inherit the bordered "description row" look from rows 7 (A:C, 3-line wrap height)
# and the date-column formatting from row 8 (D:E).
$ws.Range("A7:C7").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$ws.Range("D8:E8").Copy()
$ws.Range("D9:E9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(9).RowHeight = 43.2

$ws.Cells.Item(9, 1).Value = "Erik Meurrens"
$ws.Cells.Item(9, 2).Value = "Garage networking"
$ws.Cells.Item(9, 3).Value = "Investigating solutions to solving WiFi deadzone issue within garage. Looking into hardware that can be used to act as a network access point for the UF network."
$ws.Cells.Item(9, 4).Value = 45676
$ws.Cells.Item(9, 5).Value = 45681

# --- Row 10: Erik Meurrens / RPi configuration script ---
$ws.Range("A6:B6").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(10, 1).Value = "Erik Meurrens"
$ws.Cells.Item(10, 2).Value = "RPi configuration script"

# --- Row 11: Benjamin Simonson / RPi configuration script ---
$ws.Range("A6:B6").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 1).Value = "Benjamin Simonson"
$ws.Cells.Item(11, 2).Value = "RPi configuration script"

# Match the final cursor position left behind by the edit.
$ws.Range("A12").Select()
